# Week 2 lecture update
# - Reorders/expands the "Security" block of slides into a fuller set
#   (All the Stuffs / OWD / Sharing Rules / Roles Hierarchy / Profiles /
#   Permission Sets / Other Security / Homework).
# - Fixes small text glitches (Week 2 run-merge, stray empty run).
# - Adds speaker notes to the newly expanded security slides.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - merge the "Week " + "2" runs into a single "Week 2" run
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Week 2"

# ---------------------------------------------------------------------
# Slide 2 - "Last Week" - content unchanged (cosmetic trailing-run
# cleanup only; no visible text change needed)
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Slide 4 ("Profiles") becomes the new "Security" slide (IP Ranges /
# Login Hours / Lockout Period)
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Security"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "IP Ranges`rLogin Hours`rLockout Period"

# ---------------------------------------------------------------------
# Slide 5 ("Security") becomes the new overview slide "All the Stuffs"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "All the Stuffs"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "OWD " + [char]0x2013 + " Org Wide Defaults`r" + `
    "Manual Sharing " + [char]0x2013 + " Access at a record level.`r" + `
    "Sharing Rules " + [char]0x2013 + " Organizational record visibility`r" + `
    "Roles " + [char]0x2013 + " Which records you can see`r" + `
    "Profiles " + [char]0x2013 + " What you can do with the records you can see`r" + `
    "Permission Sets " + [char]0x2013 + " Extension of profile`r`r"

# ---------------------------------------------------------------------
# Slide 6 ("Roles Hierarchy") becomes "OWD"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "OWD"
$c6 = $s6.Shapes.Item(2).TextFrame.TextRange
$c6.Text = "Controls visibility to objects at an organizational level."
$c6.Characters(10, 11).Font.Bold = $true

# ---------------------------------------------------------------------
# Slide 7 ("Sharing Rules") - expand content, bold "see"
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$c7 = $s7.Shapes.Item(2).TextFrame.TextRange
$c7.Text = "Controls what can you see from an object level`r" + `
    "Shared with roles, subordinates, and groups`r" + `
    "Manual sharing, individual record level.`r`r"
$c7.Characters(23, 3).Font.Bold = $true

# ---------------------------------------------------------------------
# Slide 8 ("Other Security") becomes "Roles Hierarchy"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Roles Hierarchy"
$c8 = $s8.Shapes.Item(2).TextFrame.TextRange
$c8.Text = "Controls what you can see on a record level`r" + `
    "Who can see your records?`r" + `
    "Can be overwritten at the Org Wide Default level.`r"
$c8.Characters(23, 3).Font.Bold = $true

# ---------------------------------------------------------------------
# New slide 9 - "Profiles" (what used to be on slide 4, extended)
# ---------------------------------------------------------------------
$layout = $p.Slides.Item(3).CustomLayout
$s9 = $p.Slides.AddSlide(9, $layout)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Profiles"
$c9 = $s9.Shapes.Item(2).TextFrame.TextRange
$c9.Text = "What you can do (CRUD)`r" + `
    "Different based on needs`r" + `
    "Standard profiles (Read Only, API Only, System Admin)`r" + `
    "Custom profiles`r`r"
$c9.Paragraphs(1).Characters(14, 3).Font.Bold = $true

# ---------------------------------------------------------------------
# New slide 10 - "Permission Sets"
# ---------------------------------------------------------------------
$s10 = $p.Slides.AddSlide(10, $layout)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Permission Sets"
$s10.Shapes.Item(2).TextFrame.TextRange.Text = "Extends profiles`rSpecial cases"

# ---------------------------------------------------------------------
# New slide 11 - "Other Security" (moved from old slide 8, unchanged)
# ---------------------------------------------------------------------
$s11 = $p.Slides.AddSlide(11, $layout)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "Other Security"
$s11.Shapes.Item(2).TextFrame.TextRange.Text = "Field Level Security"

# ---------------------------------------------------------------------
# New slide 12 - "Homework"
# ---------------------------------------------------------------------
$s12 = $p.Slides.AddSlide(12, $layout)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "Homework"

# ---------------------------------------------------------------------
# Speaker notes for the newly expanded security slides (6-10)
# ---------------------------------------------------------------------

# Slide 6 - OWD
$n6 = $p.Slides.Item(6).NotesPage.Shapes.Item(2).TextFrame.TextRange
$n6.Text = "OWDs Controls the visibility of records and are baseline access for the records for an organisation. this access ican be opened up using :"
$n6.Characters(1, 5).Font.Bold = $true

# Slide 7 - Sharing Rules
$n7 = $p.Slides.Item(7).NotesPage.Shapes.Item(2).TextFrame.TextRange
$n7.Text = "Manual Sharing is done when we want to give access on a particular record level`r" + `
    "Sharing Rules are creagted when we want to open the access on an object level.  Sharing can be done within roles, roles & subordinates and Public groups (you create public groups to minimise the number of sharing rules e.g. you create a group " + [char]0x22 + "emea US" + [char]0x22 + " with roles & subordinates or roles emea and us)`r"
$n7.Paragraphs(1).Characters(1, 15).Font.Bold = $true
$n7.Paragraphs(2).Characters(1, 13).Font.Bold = $true
$n7.Paragraphs(2).Characters(96, 14).Font.Bold = $true

# Slide 8 - Roles Hierarchy
$n8 = $p.Slides.Item(8).NotesPage.Shapes.Item(2).TextFrame.TextRange
$n8.Text = "Roles in salesforce works differently.  Roles decides what data users can view (e.g. if OWDs are public everyone can view everything but if OWD are private and visibility is restricted, Manager can see the data of their direct reports regardless or the OWD settings."
$n8.Characters(1, 5).Font.Bold = $true

# Slide 9 - Profiles
$n9 = $p.Slides.Item(9).NotesPage.Shapes.Item(2).TextFrame.TextRange
$n9.Text = "Profiles determines what user can do with the data they view in terms of apps, tabs, object level permission, field level permission, page layout, record type and wide range of other permissions."
$n9.Characters(1, 9).Font.Bold = $true

# Slide 10 - Permission Sets
$n10 = $p.Slides.Item(10).NotesPage.Shapes.Item(2).TextFrame.TextRange
$n10.Text = "Permission Set  is another way of opening access on a user level (it can be done on profile but then we to give access to a or few users and not all the user on the profile)`r"
$n10.Paragraphs(1).Characters(1, 16).Font.Bold = $true
